# Apply the "added amoung and income to fake data 2021 and 2020" edit:
# Insert three new columns (EGID, EWID, STEUERBARESEINKOMMEN) among the
# existing VERMOGEN/HASEL/HASSH columns, and append a new AMOUNT column
# at the end, then populate header + row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before the old "VERMOGEN" column (N) ---
# This shifts VERMOGEN -> P, HASEL -> Q, HASSH -> R
$ws.Columns("N:O").Insert()

# --- Insert one new column before the (now shifted) "HASEL" column (Q) ---
# This shifts HASEL -> R, HASSH -> S
$ws.Columns("Q:Q").Insert()

# --- Header row (row 1) ---
# Note: new shared-string values are registered in the order they are first
# written, so set them in the same order the target file expects
# (EGID=36, EWID=37, AMOUNT=38, STEUERBARESEINKOMMEN=39).
$ws.Range("N1").Value = "EGID"
$ws.Range("O1").Value = "EWID"
$ws.Range("T1").Value = "AMOUNT"
$ws.Range("Q1").Value = "STEUERBARESEINKOMMEN"

# T1 is a brand-new header cell outside the old dimension, so it does not
# inherit the bold header formatting used by the rest of row 1 - copy it
# over from the adjacent header cell (S1 = "HASSH").
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

# --- Row 2 data ---
$ws.Range("N2").Value = 222
$ws.Range("O2").Value = 2
$ws.Range("Q2").Value = 25000
$ws.Range("T2").Value = -100

# --- Row 3 data ---
$ws.Range("N3").Value = 1231
$ws.Range("O3").Value = 122
$ws.Range("Q3").Value = 500000
$ws.Range("T3").Value = 777

# --- Row 4 data ---
$ws.Range("N4").Value = 9999999
$ws.Range("O4").Value = 999
$ws.Range("Q4").Value = 50000
$ws.Range("T4").Value = 99

# Leave the selection where the editor last left it before saving.
$null = $ws.Range("Q10").Select()
